$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone row 2's formatting down into rows 3-5 before filling in new content.
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2: Alerts - keep TCID, update description, clear Results
$ws.Range("B2").Value = "Creates an new Alert"
$ws.Range("E2").Value = ""

# Row 3: News Carousel
$ws.Range("A3").Value = "EXL_CorporateLensHomePage_NewsCarousel"
$ws.Range("B3").Value = "Creates a News Carousel"
$ws.Range("C3").Value = "N"
$ws.Range("D3").Value = "Y"
$ws.Range("F3").Value = "Sprint1"

# Row 4: Blogs
$ws.Range("A4").Value = "EXL_CorporateLensHomePage_Blogs"
$ws.Range("B4").Value = "Creates a New Blog"
$ws.Range("C4").Value = "N"
$ws.Range("D4").Value = "Y"
$ws.Range("F4").Value = "Sprint1"

# Row 5: My Documents
$ws.Range("A5").Value = "EXL_CorporateLensHomePage_MyDocuments"
$ws.Range("B5").Value = "Upload a New Document"
$ws.Range("C5").Value = "N"
$ws.Range("D5").Value = "Y"
$ws.Range("F5").Value = "Sprint1"

# Extend the data validation ranges to cover the new rows
$ws.Range("C2:D2").Validation.Delete()
$ws.Range("F2").Validation.Delete()
$ws.Range("C2:D5").Validation.Add(3, 1, 1, """Y,N""")
$ws.Range("F2:F5").Validation.Add(3, 1, 1, """Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10""")

# Move the selection to match the saved view
[void]$ws.Range("A2:A5").Select()
